$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 62500076
$ws.Range("I2").Value = 62500076
$ws.Range("K2").Value = 62500076
$ws.Range("M2").Value = -62499963

$ws.Range("H15").Value = 191.43
$ws.Range("I15").Value = 191.43
$ws.Range("K15").Value = 574.29
$ws.Range("M15").Value = -405.29

$ws.Range("H29").Value = 4500
$ws.Range("J29").Value = 4500
$ws.Range("L29").Value = 13500
$ws.Range("N29").Value = -14062

$ws.Range("H40").Value = 10205899
$ws.Range("I40").Value = 1798.9459
$ws.Range("J40").Value = 41668540
$ws.Range("K40").Value = 1798.9459
$ws.Range("L40").Value = 41668540
$ws.Range("M40").Value = -1623.9459
$ws.Range("N40").Value = -41668890

$ws.Range("H62").Value = 23812350
$ws.Range("I62").Value = 23812350
$ws.Range("K62").Value = 23812350
$ws.Range("M62").Value = -23811726

$ws.Range("H65").Value = 23812350
$ws.Range("I65").Value = 23812350
$ws.Range("K65").Value = 119061750
$ws.Range("M65").Value = -119058630

$ws.Range("H82").Value = 8500
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 8500
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 25500
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -26312

$ws.Range("H85").Value = 8500
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 8500
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 25500
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -28308

$ws.Range("H98").Value = 358.16666
$ws.Range("I98").Value = 368.6
$ws.Range("J98").Value = 306
$ws.Range("K98").Value = 368.6
$ws.Range("L98").Value = 306
$ws.Range("M98").Value = 1129.4
$ws.Range("N98").Value = -3302

$ws.Range("H122").Value = 358.16666
$ws.Range("I122").Value = 368.6
$ws.Range("J122").Value = 306
$ws.Range("K122").Value = 1105.8
$ws.Range("L122").Value = 918
$ws.Range("M122").Value = 1344.2
$ws.Range("N122").Value = -5818

$ws.Range("H132").Value = 4696880.5
$ws.Range("I132").Value = 1518.5079
$ws.Range("K132").Value = 4555.5237
$ws.Range("M132").Value = -2025.5237

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18534.803
$ws.Range("I32").Value = 14854.346
$ws.Range("J32").Value = 31186.375
$ws.Range("K32").Value = 14854.346
$ws.Range("L32").Value = 31186.375
$ws.Range("M32").Value = -14567.346
$ws.Range("N32").Value = -31760.375

$ws.Range("H45").Value = 151103.86
$ws.Range("I45").Value = 191632.19
$ws.Range("J45").Value = 2500
$ws.Range("K45").Value = 191632.19
$ws.Range("L45").Value = 2500
$ws.Range("M45").Value = -191255.19
$ws.Range("N45").Value = -3254

$ws.Range("H61").Value = 179528.67
$ws.Range("I61").Value = 4677.324
$ws.Range("J61").Value = 503003.66
$ws.Range("K61").Value = 4677.324
$ws.Range("L61").Value = 503003.66
$ws.Range("M61").Value = -4465.324
$ws.Range("N61").Value = -503427.66

$ws.Range("H97").Value = 734.2353000000001
$ws.Range("I97").Value = 511.5
$ws.Range("J97").Value = 1268.8
$ws.Range("K97").Value = 511.5
$ws.Range("L97").Value = 1268.8
$ws.Range("M97").Value = -15.5
$ws.Range("N97").Value = -2260.8

$ws.Range("H132").Value = 2798.5881
$ws.Range("I132").Value = 1737.7354
$ws.Range("J132").Value = 4920.294
$ws.Range("K132").Value = 5213.206200000001
$ws.Range("L132").Value = 14760.882
$ws.Range("M132").Value = -2683.206200000001
$ws.Range("N132").Value = -19820.882

$ws.Range("H136").Value = 179528.67
$ws.Range("I136").Value = 4677.324
$ws.Range("J136").Value = 503003.66
$ws.Range("K136").Value = 14031.972
$ws.Range("L136").Value = 1509010.98
$ws.Range("M136").Value = -11481.972
$ws.Range("N136").Value = -1514110.98

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 300720.1
$ws.Range("I107").Value = 300720.1
$ws.Range("K107").Value = 300720.1
$ws.Range("M107").Value = -298800.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9809158
$ws.Range("I31").Value = 1654.8846
$ws.Range("J31").Value = 20008962
$ws.Range("K31").Value = 1654.8846
$ws.Range("L31").Value = 20008962
$ws.Range("M31").Value = -1359.8846
$ws.Range("N31").Value = -20009552

$ws.Range("H34").Value = 9809158
$ws.Range("I34").Value = 1654.8846
$ws.Range("J34").Value = 20008962
$ws.Range("K34").Value = 1654.8846
$ws.Range("L34").Value = 20008962
$ws.Range("M34").Value = -1452.8846
$ws.Range("N34").Value = -20009366

$ws.Range("H58").Value = 5355136
$ws.Range("I58").Value = 6174391
$ws.Range("J58").Value = 1668488.1
$ws.Range("K58").Value = 6174391
$ws.Range("L58").Value = 1668488.1
$ws.Range("M58").Value = -6174188
$ws.Range("N58").Value = -1668894.1

$ws.Range("H99").Value = 2126.3
$ws.Range("I99").Value = 1390.25
$ws.Range("J99").Value = 2617
$ws.Range("K99").Value = 1390.25
$ws.Range("L99").Value = 2617
$ws.Range("M99").Value = 107.75
$ws.Range("N99").Value = -5613

$ws.Range("H126").Value = 2126.3
$ws.Range("I126").Value = 1390.25
$ws.Range("J126").Value = 2617
$ws.Range("K126").Value = 4170.75
$ws.Range("L126").Value = 7851
$ws.Range("M126").Value = -1700.75
$ws.Range("N126").Value = -12791

$ws.Range("H132").Value = 4445887
$ws.Range("I132").Value = 4879232
$ws.Range("J132").Value = 4100.25
$ws.Range("K132").Value = 14637696
$ws.Range("L132").Value = 12300.75
$ws.Range("M132").Value = -14635166
$ws.Range("N132").Value = -17360.75

$ws.Range("H136").Value = 5355136
$ws.Range("I136").Value = 6174391
$ws.Range("J136").Value = 1668488.1
$ws.Range("K136").Value = 18523173
$ws.Range("L136").Value = 5005464.300000001
$ws.Range("M136").Value = -18520623
$ws.Range("N136").Value = -5010564.300000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2273821.2
$ws.Range("J131").Value = 1202.0605
$ws.Range("L131").Value = 3606.1815
$ws.Range("N131").Value = -13686.1815

$ws.Range("H137").Value = 17741.334
$ws.Range("J137").Value = 26357.455
$ws.Range("L137").Value = 79072.36500000001
$ws.Range("N137").Value = -89272.36500000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3398.8572
$ws.Range("I102").Value = 2969.75
$ws.Range("J102").Value = 3971
$ws.Range("K102").Value = 2969.75
$ws.Range("L102").Value = 3971
$ws.Range("M102").Value = -1347.75
$ws.Range("N102").Value = -7215

$ws.Range("H122").Value = 113737384
$ws.Range("I122").Value = 177470640
$ws.Range("K122").Value = 532411920
$ws.Range("M122").Value = -532409470

$ws.Range("H132").Value = 5052273
$ws.Range("I132").Value = 5748552
$ws.Range("J132").Value = 4249.75
$ws.Range("K132").Value = 17245656
$ws.Range("L132").Value = 12749.25
$ws.Range("M132").Value = -17243126
$ws.Range("N132").Value = -17809.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2992.5454
$ws.Range("I7").Value = 1956.4445
$ws.Range("J7").Value = 7655
$ws.Range("K7").Value = 1956.4445
$ws.Range("L7").Value = 7655
$ws.Range("M7").Value = -1844.4445
$ws.Range("N7").Value = -7879

$ws.Range("H126").Value = 2992.5454
$ws.Range("I126").Value = 1956.4445
$ws.Range("J126").Value = 7655
$ws.Range("K126").Value = 5869.333500000001
$ws.Range("L126").Value = 22965
$ws.Range("M126").Value = -3399.333500000001
$ws.Range("N126").Value = -27905

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1168.5385
$ws.Range("I126").Value = 1157.6
$ws.Range("J126").Value = 1205
$ws.Range("K126").Value = 3472.8
$ws.Range("L126").Value = 3615
$ws.Range("M126").Value = -1002.8
$ws.Range("N126").Value = -8555

$ws.Range("H132").Value = 1820.55
$ws.Range("I132").Value = 771.04346
$ws.Range("J132").Value = 3240.4707
$ws.Range("K132").Value = 2313.13038
$ws.Range("L132").Value = 9721.4121
$ws.Range("M132").Value = 216.8696199999999
$ws.Range("N132").Value = -14781.4121

$ws.Range("H136").Value = 3088602
$ws.Range("I136").Value = 2069.4736
$ws.Range("J136").Value = 10419117
$ws.Range("K136").Value = 6208.4208
$ws.Range("L136").Value = 31257351
$ws.Range("M136").Value = -3658.4208
$ws.Range("N136").Value = -31262451
